$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("model_6_5_24",0.3353703829575784,0.1723236691864223,0.286546990682506,0.2278705744312793,0.7355490922927856,0.8305639624595642,0.7144161462783813,0.7759062051773071),
    @("model_6_5_23",0.3354439741309539,0.1734385344459174,0.2867116985065469,0.2285372800963296,0.735467791557312,0.8294452428817749,0.7142511606216431,0.7752362489700317),
    @("model_6_5_10",0.3355978892247333,0.180900176682678,0.2852391853402723,0.2317914881687008,0.7352973222732544,0.8219575881958008,0.7157256603240967,0.7719660997390747),
    @("model_6_5_16",0.3356623481599501,0.1780290613492953,0.2862204061003548,0.2307336396793654,0.7352259755134583,0.8248386979103088,0.7147431373596191,0.773029088973999),
    @("model_6_5_17",0.3357136845317279,0.1779445222907166,0.2865070016727838,0.2308235130114579,0.7351691126823425,0.824923574924469,0.7144562005996704,0.7729387879371643),
    @("model_6_5_11",0.33572671149866,0.1809835590513745,0.2856091793408669,0.2320091218443801,0.7351547479629517,0.8218739032745361,0.7153552770614624,0.7717474102973938),
    @("model_6_5_22",0.3357437285213457,0.1756800449111081,0.2876110376688084,0.2301440566335434,0.7351359724998474,0.8271958827972412,0.7133506536483765,0.7736216187477112),
    @("model_6_5_21",0.335760580704252,0.1766684821441411,0.2869958160824779,0.2303778028670567,0.7351173162460327,0.8262040019035339,0.7139666676521301,0.7733866572380066),
    @("model_6_5_8",0.3357663582638823,0.1822943824008885,0.2855854130472314,0.2326910070594308,0.7351108193397522,0.820558488368988,0.715378999710083,0.7710621953010559),
    @("model_6_5_9",0.3357695337580935,0.1821558736373831,0.285588008820734,0.232618879502688,0.735107421875,0.8206974864006042,0.7153764367103577,0.7711346745491028),
    @("model_6_5_7",0.3357759438952828,0.182699101352487,0.2856806886785053,0.2329493763937521,0.7351003289222717,0.8201524019241333,0.7152836322784424,0.7708025574684143),
    @("model_6_5_13",0.3357816076296807,0.1807690312005926,0.2858246372246326,0.2319966677788491,0.7350940704345703,0.8220891356468201,0.7151395082473755,0.7717599272727966),
    @("model_6_5_20",0.3357858268454584,0.1769271906691872,0.2870330599048813,0.2305323107708583,0.7350893616676331,0.8259443640708923,0.7139294147491455,0.7732314467430115),
    @("model_6_5_19",0.3357870202439798,0.1771751652631153,0.2869753402933568,0.2306363512158406,0.735088050365448,0.8256955742835999,0.7139872312545776,0.7731269001960754),
    @("model_6_5_18",0.3358187766295931,0.1775321266712424,0.2870851861368436,0.2308764884204964,0.7350529432296753,0.8253373503684998,0.7138772010803223,0.7728855609893799),
    @("model_6_5_12",0.3358204069446715,0.1810714297274236,0.2859404481620518,0.2322109695685209,0.7350510954856873,0.8217856884002686,0.7150235176086426,0.7715445756912231),
    @("model_6_5_15",0.335887203879186,0.1806918506283303,0.2862224761534649,0.2321422281524419,0.7349771857261658,0.8221666812896729,0.7147411108016968,0.771613597869873),
    @("model_6_5_14",0.3358880269431532,0.1807168445588673,0.2862210144659229,0.2321549912511044,0.7349762916564941,0.8221415281295776,0.7147425413131714,0.7716007828712463),
    @("model_6_5_0",0.335956188912967,0.2200996527014069,0.302626316535718,0.2606685087307518,0.7349008321762085,0.7826213836669922,0.6983150839805603,0.7429478168487549),
    @("model_6_5_5",0.3360093213602412,0.1866228201690906,0.2859805505309602,0.2351646082839467,0.7348420023918152,0.816214919090271,0.714983344078064,0.7685765027999878),
    @("model_6_5_6",0.3360494293516049,0.1854522602945878,0.2864963282162041,0.2347877084958594,0.7347975969314575,0.8173896074295044,0.7144668102264404,0.7689552307128906),
    @("model_6_5_3",0.3362890096266673,0.2025401667752621,0.2860085317807763,0.2435925931780798,0.7345324754714966,0.8002421259880066,0.7149553298950195,0.7601072192192078),
    @("model_6_5_4",0.3363109760579339,0.2025964571815095,0.2860067610307773,0.2436216985589176,0.7345081567764282,0.8001855611801147,0.714957058429718,0.760077953338623),
    @("model_6_5_2",0.3364755146737628,0.2036063053894963,0.2865130388891459,0.2443928467154106,0.7343260645866394,0.7991722226142883,0.7144501209259033,0.759303092956543),
    @("model_6_5_1",0.3364814299372867,0.2037644400486615,0.289184319402312,0.2457293048741627,0.7343195676803589,0.7990134358406067,0.7117752432823181,0.7579600811004639),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
